$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed the new rows' formatting (column D date style) from the last existing
# row so we reuse the existing style index instead of minting a new one.
$ws.Range("D70").Copy($ws.Range("D71"))
$ws.Range("D70").Copy($ws.Range("D72"))

# Two new work-log entries appended below the existing data (rows 71-72)
$ws.Cells.Item(71, 2).Value = "Psaní - rešerše + screenshoty, řešení vlastních enumerací"
$ws.Cells.Item(71, 3).Value = 4
$ws.Cells.Item(71, 4).Value = (Get-Date -Year 2012 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item(72, 2).Value = "Psaní - bezpečnost, řešení citací, řešení a  tvorba příkazů pro vkládání zdrojového kódu C#"
$ws.Cells.Item(72, 3).Value = 7
$ws.Cells.Item(72, 4).Value = (Get-Date -Year 2012 -Month 4 -Day 9 -Hour 0 -Minute 0 -Second 0)

# Leave the selection where the user ended up after entering the new rows
$ws.Range("B73").Select()

$wb.Save()
